$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 17) to the BIIBNoun sheet, mirroring the
# format already used by the existing data rows (2-16).
$srcRow = 16
$row = 17

# Copy formatting (style) from the last existing data row so the new
# row's date cell keeps the same number format / style index.
$ws.Range("A" + $srcRow + ":N" + $srcRow).Copy() | Out-Null
$ws.Range("A" + $row + ":N" + $row).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 42622.888460648152
$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 58
$ws.Cells.Item($row, 4).Value = 37
$ws.Cells.Item($row, 5).Value = 58
$ws.Cells.Item($row, 6).Value = 22
$ws.Cells.Item($row, 7).Value = 15507
$ws.Cells.Item($row, 8).Value = 12785
$ws.Cells.Item($row, 9).Value = 2082
$ws.Cells.Item($row, 10).Value = 293
$ws.Cells.Item($row, 11).Value = 187
$ws.Cells.Item($row, 12).Value = 44
$ws.Cells.Item($row, 13).Value = 13
$ws.Cells.Item($row, 14).Value = "Noun"

$wb.Save()
